$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added to the dataset. In the source sheet
# this corresponds to inserting a new row at row 606, which pushes all the
# following rows (old 606-648) down by one (new 607-649), and filling the
# newly inserted row 606 with the new record's data.
$ws.Rows("606:606").Insert()

$ws.Range("A606").Value = 11
$ws.Range("B606").Value = "Vega Monumental Concepción"
$ws.Range("C606").Value = "Bíobío"
$ws.Range("D606").Value = 44931
$ws.Range("E606").Value = 8
$ws.Range("F606").Value = 100112004
$ws.Range("G606").Value = "Cebolla"
$ws.Range("H606").Value = "Sin especificar"
$ws.Range("I606").Value = "1a (cosecha)"
$ws.Range("J606").Value = 450
$ws.Range("K606").Value = 8000
$ws.Range("L606").Value = 8500
$ws.Range("M606").Value = 8278
$ws.Range("N606").Value = "$/malla 18 kilos"
$ws.Range("O606").Value = "Región Metropolitana"
$ws.Range("P606").Value = 460
$ws.Range("Q606").Value = 18
$ws.Range("R606").Value = "Hortaliza"
